# Fix the run-together hex-color summary line in the "Highlights" section so
# that each color code gets its own paragraph along with its count, e.g.
#   "#fb5b89#69aff0#7cc867#f9cd59"
# becomes four separate paragraphs:
#   "#fb5b89: 43"
#   "#69aff0: 18"
#   "#7cc867: 16"
#   "#f9cd59: 20"

$d = $word.ActiveDocument

$oldText = "#fb5b89#69aff0#7cc867#f9cd59"
$newText = "#fb5b89: 43" + [char]13 + "#69aff0: 18" + [char]13 + "#7cc867: 16" + [char]13 + "#f9cd59: 20"

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq $oldText) {
        $p.Range.Text = $newText
        break
    }
}
